# Generate Report for Handoff
#
# The localization status report has moved from "In Translation" to
# "Ready for handoff": update the Status text and the associated
# timestamps on all three sheets, then widen the Status column(s) so the
# longer status string still fits (mirrors Excel's column autofit).

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- Overview sheet --------------------------------------------------
# Column E = zh-cn status, Column F = de-de status,
# Column G = Latest HO Xliff Generate Date
$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
$overview.Range("G2").Value = "2016-08-20 15:00:47"

# --- zh-cn sheet -------------------------------------------------------
# Column C = Status, Column H = Latest Handoff Datetime
$zhcn.Range("C2").Value = "Ready for handoff"
$zhcn.Range("H2").Value = "2016-08-20 15:00:43"

# --- de-de sheet -------------------------------------------------------
# Column C = Status, Column H = Latest Handoff Datetime
$dede.Range("C2").Value = "Ready for handoff"
$dede.Range("H2").Value = "2016-08-20 15:00:47"

# --- Widen the Status columns to fit "Ready for handoff" ---------------
$overview.Columns.Item(5).ColumnWidth = 16.333333333333332
$overview.Columns.Item(6).ColumnWidth = 16.333333333333332
$zhcn.Columns.Item(3).ColumnWidth = 16.333333333333332
$dede.Columns.Item(3).ColumnWidth = 16.333333333333332
